$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Cells changing from shared-text placeholder to a real number (need target number format first,
#     so Excel reuses the existing numeric style instead of minting a new one) ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 3
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -66.666666666666
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 3
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H15").Value = -33.333333333333
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 3
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -66.666666666666
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 3
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H26").Value = 66.666666666666

# --- Cells changing from a real number back to the shared-text placeholder ("0" / "***.*") ---
#     Copy from a donor cell that already holds that exact placeholder + style so the shared string
#     and style index line up with the rest of the sheet. ---
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("C27"))

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("L14").Value = -75
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -42.307692307692
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -40.909090909090
$ws.Range("I16").Value = 140
$ws.Range("J16").Value = 158
$ws.Range("K16").Value = -11.392405063291
$ws.Range("L16").Value = 8.527131782945
$ws.Range("M16").Value = -19.540229885057
$ws.Range("N16").Value = -76.068376068376
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -27.586206896551
$ws.Range("I17").Value = 245
$ws.Range("J17").Value = 245
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12.385321100917
$ws.Range("M17").Value = 160.63829787234
$ws.Range("N17").Value = -14.035087719298
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -23.076923076923
$ws.Range("I18").Value = 98
$ws.Range("K18").Value = -2
$ws.Range("L18").Value = 60.655737704918
$ws.Range("M18").Value = -51.485148514851
$ws.Range("N18").Value = -87.901234567901
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 116.666666666667
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 40
$ws.Range("I19").Value = 432
$ws.Range("J19").Value = 466
$ws.Range("K19").Value = -7.296137339055
$ws.Range("L19").Value = 33.333333333333
$ws.Range("M19").Value = 86.206896551724
$ws.Range("N19").Value = 6.666666666666
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 21
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 172
$ws.Range("J20").Value = 153
$ws.Range("K20").Value = 12.418300653594
$ws.Range("L20").Value = 43.333333333333
$ws.Range("M20").Value = -13.131313131313
$ws.Range("N20").Value = -92.293906810035
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 29.166666666666
$ws.Range("F21").Value = 116
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = 1.754385964912
$ws.Range("I21").Value = 1103
$ws.Range("J21").Value = 1144
$ws.Range("K21").Value = -3.583916083916
$ws.Range("L21").Value = 27.367205542725
$ws.Range("M21").Value = 20.021762785636
$ws.Range("N21").Value = -74.696031199816
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 17
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 41.666666666666
$ws.Range("M22").Value = -5.555555555555
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -35.606060606060
$ws.Range("I24").Value = 889
$ws.Range("J24").Value = 990
$ws.Range("K24").Value = -10.202020202020
$ws.Range("L24").Value = 46.457990115321
$ws.Range("M24").Value = 100.677200902935
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -24.528301886792
$ws.Range("I25").Value = 382
$ws.Range("J25").Value = 375
$ws.Range("K25").Value = 1.866666666666
$ws.Range("L25").Value = 20.886075949367
$ws.Range("M25").Value = 4.945054945054
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 28
$ws.Range("K26").Value = -14.285714285714
$ws.Range("L26").Value = 60
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -5.128205128205
$ws.Range("L28").Value = -66.666666666666
$ws.Range("L29").Value = -72.727272727272
